# Automatische test-sync: 2025-07-29 21:35:50
#
# Adds a 4th test-mail row to the "Logs" sheet, extends its conditional
# formatting ranges to cover the new row, adds the matching roll-up row to
# the "Dashboard" sheet, and extends the bar chart's category/value series
# references so the new "Bestelling / Levering" bucket is plotted.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 6 with the new test e-mail.
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A6").Value = "Wil je 100 stuks M5-bouten bestellen?"
$logs.Range("B6").Value = "mailmind.test@zohomail.eu"
$logs.Range("C6").Value = "Testmail #4: Wil je 100 stuks M5-bouten bestellen?"
$logs.Range("D6").Value = "Bestelling / Levering"
$logs.Range("E6").Value = "Geachte afzender,`nDank u voor uw e-mail. Helaas kunnen we geen bestellingen plaatsen via deze e-mail. Gelieve onze website te bezoeken en het bestelproces te doorlopen om uw M5-bouten te bestellen.`nMet vriendelijke groet,`n[E-mailassistent]"
$logs.Range("F6").Value = "2025-07-29 21:35:19"
$logs.Range("G6").Value = "Ja"
$logs.Range("H6").Value = "Nee"
$logs.Range("I6").Value = "Ja"
$logs.Range("J6").Value = "Nee"

# The multi-line Antwoord text makes Excel auto-grow the row; AutoFit keeps
# the row height "default" (no custom ht) instead of leaving a stray
# explicit row height behind, matching the other (unstyled) rows.
$logs.Rows.Item(6).AutoFit()

# Extend the conditional-formatting ranges (D/G/H/I/J, rows 2:5) to 2:6 so
# the new row is covered. Each FormatConditions collection shares one
# <conditionalFormatting sqref="..."> block, so touching rule 1 updates
# the whole group.
$logs.Range("D2:D5").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D6"))
$logs.Range("G2:G5").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G6"))
$logs.Range("H2:H5").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H6"))
$logs.Range("I2:I5").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I6"))
$logs.Range("J2:J5").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J6"))

# ---------------------------------------------------------------------
# 2. Dashboard sheet: append the matching category roll-up row.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A4").Value = "Bestelling / Levering"
$dash.Range("B4").Value = 1

# ---------------------------------------------------------------------
# 3. Chart: extend the category/value series references from row 3 to
#    row 4 so the new Dashboard row is plotted.
# ---------------------------------------------------------------------
$chart = $dash.ChartObjects().Item(1).Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$4,'Dashboard'!`$B`$2:`$B`$4,1)"
